$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # users
$ws2 = $wb.Worksheets.Item(2)   # products
$ws3 = $wb.Worksheets.Item(3)   # orders

# ---------------------------------------------------------------------------
# New "create_date" / "create_time" columns on the "orders" sheet (F & G),
# added to support DATE_FORMAT validation against SQLite / PostgreSQL dumps.
# ---------------------------------------------------------------------------

$headerStyleSource = $ws3.Range("E1")

# Header cells (F1/G1) - written in this exact order so the generated
# sharedStrings table lines up with the source fixture.
$ws3.Range("F1").Value = "create_date"
$headerStyleSource.Copy()
$ws3.Range("F1").PasteSpecial(-4122)

# create_date values (column F) - entered out of row order to mirror the
# fixture's original insertion sequence; kept as literal text (not dates).
$ws3.Range("F2:F6").NumberFormat = "@"
$ws3.Range("F2").Value = "03/04/2004"
$ws3.Range("F5").Value = "23/04/2007"
$ws3.Range("F6").Value = "02/31/2008"
$ws3.Range("F4").Value = "03/24/2006"
$ws3.Range("F3").Value = "11/07/2005"
$ws3.Range("F2:F6").ClearFormats()

$ws3.Range("G1").Value = "create_time"
$headerStyleSource.Copy()
$ws3.Range("G1").PasteSpecial(-4122)

# create_time values (column G) - kept as literal text via an explicit
# Text number format (covers the intentionally-invalid 02/31/2008, 12:33:78
# and 25:33:48 samples used to exercise DATE_FORMAT validation).
$ws3.Range("G2:G6").NumberFormat = "@"
$ws3.Range("G2").Value = "12:33:44"
$ws3.Range("G4").Value = "12:33:46"
$ws3.Range("G5").Value = "12:33:47"
$ws3.Range("G3").Value = "12:33:78"
$ws3.Range("G6").Value = "25:33:48"

# Column widths for the new columns.
$ws3.Columns.Item(4).ColumnWidth = 11.6
$ws3.Columns.Item(5).ColumnWidth = 12.6
$ws3.Columns.Item(6).ColumnWidth = 12.6
$ws3.Columns.Item(7).ColumnWidth = 12.6

# Print setup for the "orders" sheet.
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping (matches the recorded UI state).
# ---------------------------------------------------------------------------
$ws1.Range("F3").Select()
$ws2.Range("E6").Select()
$ws3.Range("G4").Select()
$ws3.Activate()
